$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D price cells to Text format so numeric-looking strings
# (e.g. "1.00", "54.853.96") are preserved exactly as text, matching the
# original inline-string cell type instead of being coerced to numbers.
$dCells = @("D2","D3","D5","D6","D7","D8","D9","D10","D12","D13","D14","D15","D16","D18","D19","D21","D22","D23","D24","D25","D27","D28","D29","D30","D31","D32","D34","D35","D36","D38","D39","D40","D41","D42","D43","D44","D45","D46","D49","D50","D51")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '54.853.96'
$ws.Range('E2').Value = '  +1.02%  '
$ws.Range('D3').Value = '2.292.63'
$ws.Range('E3').Value = '  +0.58%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '507.44'
$ws.Range('E5').Value = '  +0.79%  '
$ws.Range('D6').Value = '129.75'
$ws.Range('E6').Value = '  +0.37%  '
$ws.Range('D7').Value = '0.995'
$ws.Range('E7').Value = '  -0.38%  '
$ws.Range('D8').Value = '0.531'
$ws.Range('E8').Value = '  +0.37%  '
$ws.Range('D9').Value = '2.316.75'
$ws.Range('E9').Value = '  +1.20%  '
$ws.Range('D10').Value = '0.0973'
$ws.Range('E10').Value = '  +1.91%  '
$ws.Range('E11').Value = '  +2.03%  '
$ws.Range('D12').Value = '0.340'
$ws.Range('E12').Value = '  +2.01%  '
$ws.Range('D13').Value = '4.99'
$ws.Range('E13').Value = '  +5.50%  '
$ws.Range('D14').Value = '23.90'
$ws.Range('E14').Value = '  +4.51%  '
$ws.Range('D15').Value = '2.702.67'
$ws.Range('E15').Value = '  +0.53%  '
$ws.Range('D16').Value = '54.853.01'
$ws.Range('E16').Value = '  +1.09%  '
$ws.Range('E17').Value = '  +1.76%  '
$ws.Range('D18').Value = '2.405.19'
$ws.Range('E18').Value = '  +5.34%  '
$ws.Range('D19').Value = '10.65'
$ws.Range('E19').Value = '  +3.90%  '
$ws.Range('E20').Value = '  +1.79%  '
$ws.Range('D21').Value = '6.70'
$ws.Range('E21').Value = '  +4.31%  '
$ws.Range('D22').Value = '309.82'
$ws.Range('E22').Value = '  +1.55%  '
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').Value = '  +0.39%  '
$ws.Range('D24').Value = '60.53'
$ws.Range('E24').Value = '  -2.33%  '
$ws.Range('D25').Value = '0.993'
$ws.Range('E25').Value = '  -0.61%  '
$ws.Range('E26').Value = '  +0.44%  '
$ws.Range('D27').Value = '7.52'
$ws.Range('E27').Value = '  +2.58%  '
$ws.Range('D28').Value = '171.92'
$ws.Range('E28').Value = '  -1.51%  '
$ws.Range('D29').Value = '6.16'
$ws.Range('E29').Value = '  +2.53%  '
$ws.Range('B30').Value = 'PEPE'
$ws.Range('C30').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D30').Value = '0.0₃0711'
$ws.Range('E30').Value = '  +3.05%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').Value = '1.64'
$ws.Range('E31').Value = '  +0.64%  '
$ws.Range('B32').Value = 'Fetch.AI'
$ws.Range('C32').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D32').Value = '1.15'
$ws.Range('E32').Value = '  +6.39%  '
$ws.Range('E33').Value = '  -0.02%  '
$ws.Range('D34').Value = '18.06'
$ws.Range('E34').Value = '  +1.52%  '
$ws.Range('D35').Value = '0.993'
$ws.Range('E35').Value = '  -0.40%  '
$ws.Range('D36').Value = '0.915'
$ws.Range('E36').Value = '  -3.56%  '
$ws.Range('E37').Value = '  +2.18%  '
$ws.Range('D38').Value = '3.87'
$ws.Range('E38').Value = '  +3.39%  '
$ws.Range('D39').Value = '36.74'
$ws.Range('E39').Value = '  +1.84%  '
$ws.Range('D40').Value = '0.378'
$ws.Range('E40').Value = '  +1.10%  '
$ws.Range('D41').Value = '1.44'
$ws.Range('E41').Value = '  +2.11%  '
$ws.Range('D42').Value = '133.98'
$ws.Range('E42').Value = '  +6.92%  '
$ws.Range('D43').Value = '3.44'
$ws.Range('E43').Value = '  +1.57%  '
$ws.Range('D44').Value = '4.98'
$ws.Range('E44').Value = '  +3.51%  '
$ws.Range('D45').Value = '254.32'
$ws.Range('E45').Value = '  +5.90%  '
$ws.Range('D46').Value = '0.0505'
$ws.Range('E46').Value = '  +1.82%  '
$ws.Range('E47').Value = '  +2.17%  '
$ws.Range('E48').Value = '  +1.21%  '
$ws.Range('D49').Value = '0.378'
$ws.Range('E49').Value = '  +1.46%  '
$ws.Range('D50').Value = '0.0210'
$ws.Range('E50').Value = '  +1.90%  '
$ws.Range('D51').Value = '10.81'
$ws.Range('E51').Value = '  +0.38%  '
